$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.25
$ws.Range("I2").Value = 3.25
$ws.Range("J2").Value = 1.07
$ws.Range("K2").Value = 9
$ws.Range("AE2").Value = 9.5
$ws.Range("AG2").Value = 12
$ws.Range("N4").Value = 2.5
$ws.Range("O4").Value = 1.5
$ws.Range("X4").Value = 21
$ws.Range("G5").Value = 1.47
$ws.Range("H5").Value = 3.8
$ws.Range("I5").Value = 6.8
$ws.Range("L5").Value = 1.36
$ws.Range("M5").Value = 2.67
$ws.Range("N5").Value = 2.05
$ws.Range("O5").Value = 1.6
$ws.Range("P5").Value = 1.42
$ws.Range("Q5").Value = 2.47
$ws.Range("R5").Value = 2.22
$ws.Range("S5").Value = 1.52
$ws.Range("T5").Value = 5.2
$ws.Range("U5").Value = 5.8
$ws.Range("V5").Value = 8.75
$ws.Range("W5").Value = 9.25
$ws.Range("X5").Value = 14
$ws.Range("Y5").Value = 40
$ws.Range("Z5").Value = 8
$ws.Range("AA5").Value = 7.8
$ws.Range("AB5").Value = 25
$ws.Range("AE5").Value = 14
$ws.Range("AF5").Value = 40
$ws.Range("G6").Value = 2.02
$ws.Range("H6").Value = 3.15
$ws.Range("I6").Value = 3.55
$ws.Range("M6").Value = 2.35
$ws.Range("N6").Value = 2.32
$ws.Range("Q6").Value = 2.27
$ws.Range("R6").Value = 2.07
$ws.Range("T6").Value = 5.6
$ws.Range("U6").Value = 8.25
$ws.Range("V6").Value = 9.25
$ws.Range("W6").Value = 17.5
$ws.Range("X6").Value = 20
$ws.Range("AA6").Value = 6.3
$ws.Range("AB6").Value = 20
$ws.Range("AE6").Value = 8
$ws.Range("AF6").Value = 17
$ws.Range("AG6").Value = 13.5
$ws.Range("AH6").Value = 50
$ws.Range("AI6").Value = 40
$ws.Range("AJ6").Value = 60
$ws.Range("I9").Value = 3.3
$ws.Range("L9").Value = 1.2
$ws.Range("T9").Value = 9.5
$ws.Range("X9").Value = 14
$ws.Range("Z9").Value = 13
$ws.Range("AI9").Value = 26
$ws.Range("AJ9").Value = 29
$ws.Range("G13").Value = 2.15
$ws.Range("H13").Value = 3.05
$ws.Range("I13").Value = 3.15
$ws.Range("N13").Value = 2.07
$ws.Range("O13").Value = 1.6
$ws.Range("P13").Value = 1.42
$ws.Range("Q13").Value = 2.37
$ws.Range("T13").Value = 5.6
$ws.Range("U13").Value = 8
$ws.Range("V13").Value = 7.6
$ws.Range("W13").Value = 16.5
$ws.Range("X13").Value = 15.5
$ws.Range("Y13").Value = 26
$ws.Range("Z13").Value = 7.9
$ws.Range("AA13").Value = 5.3
$ws.Range("AC13").Value = 60
$ws.Range("AD13").Value = 450
$ws.Range("AE13").Value = 7.3
$ws.Range("AF13").Value = 13
$ws.Range("AG13").Value = 9.5
$ws.Range("AH13").Value = 32
$ws.Range("AI13").Value = 23
$ws.Range("AJ13").Value = 30
$ws.Range("G14").Value = 3.45
$ws.Range("H14").Value = 3.25
$ws.Range("I14").Value = 1.95
$ws.Range("N14").Value = 2.12
$ws.Range("O14").Value = 1.57
$ws.Range("Q14").Value = 2.32
$ws.Range("T14").Value = 7.4
$ws.Range("U14").Value = 14
$ws.Range("V14").Value = 10.5
$ws.Range("W14").Value = 37
$ws.Range("X14").Value = 28
$ws.Range("Y14").Value = 37
$ws.Range("Z14").Value = 7.8
$ws.Range("AA14").Value = 5.6
$ws.Range("AE14").Value = 5.3
$ws.Range("AF14").Value = 7.1
$ws.Range("AG14").Value = 7.5
$ws.Range("AH14").Value = 13.5
$ws.Range("AI14").Value = 14
$ws.Range("AJ14").Value = 27
$ws.Range("J15").Value = 1.1
$ws.Range("K15").Value = 7
$ws.Range("N15").Value = 2.35
$ws.Range("O15").Value = 1.57
$ws.Range("T15").Value = 8.5
$ws.Range("Z15").Value = 7
$ws.Range("AD15").Value = 401
$ws.Range("N16").Value = 2.2
$ws.Range("O16").Value = 1.65
$ws.Range("G17").Value = 2.45
$ws.Range("I17").Value = 2.77
$ws.Range("K17").Value = 6.4
$ws.Range("N17").Value = 2.1
$ws.Range("S17").Value = 1.88
$ws.Range("T17").Value = 7.5
$ws.Range("U17").Value = 12
$ws.Range("V17").Value = 9.5
$ws.Range("W17").Value = 27
$ws.Range("X17").Value = 21
$ws.Range("Z17").Value = 6.4
$ws.Range("AD17").Value = 600
$ws.Range("AE17").Value = 7.8
$ws.Range("AF17").Value = 13.5
$ws.Range("AG17").Value = 10.5
$ws.Range("AH17").Value = 35
$ws.Range("AI17").Value = 26
$ws.Range("AJ17").Value = 37
$ws.Range("J18").Value = 1.08
$ws.Range("K18").Value = 6.4
$ws.Range("L18").Value = 1.39
$ws.Range("M18").Value = 2.77
$ws.Range("N18").Value = 2.15
$ws.Range("O18").Value = 1.62
$ws.Range("Q18").Value = 2.42
$ws.Range("R18").Value = 1.93
$ws.Range("S18").Value = 1.78
$ws.Range("Y18").Value = 32
$ws.Range("Z18").Value = 6.4
$ws.Range("AB18").Value = 16.5
$ws.Range("AC18").Value = 90
$ws.Range("AD18").Value = 800
$ws.Range("AE18").Value = 9.5
$ws.Range("AF18").Value = 19.5
$ws.Range("AG18").Value = 13
$ws.Range("AI18").Value = 37
$ws.Range("AJ18").Value = 50
$ws.Range("T23").Value = 8.75
$ws.Range("U23").Value = 11.5
$ws.Range("Y23").Value = 17.5
$ws.Range("Z23").Value = 12
$ws.Range("AE23").Value = 9
$ws.Range("AF23").Value = 12.5
$ws.Range("AH23").Value = 25
$ws.Range("AI23").Value = 17
$ws.Range("AJ23").Value = 20
$ws.Range("G24").Value = 2.05
$ws.Range("H24").Value = 3.35
$ws.Range("I24").Value = 3.1
$ws.Range("T24").Value = 6.8
$ws.Range("U24").Value = 8.75
$ws.Range("V24").Value = 7.4
$ws.Range("W24").Value = 15.5
$ws.Range("X24").Value = 13
$ws.Range("Z24").Value = 10.75
$ws.Range("AA24").Value = 5.8
$ws.Range("AB24").Value = 11
$ws.Range("AC24").Value = 40
$ws.Range("AD24").Value = 250
$ws.Range("AE24").Value = 9
$ws.Range("AF24").Value = 14
$ws.Range("AG24").Value = 9.25
$ws.Range("AH24").Value = 32
$ws.Range("AI24").Value = 21
$ws.Range("AJ24").Value = 24
$ws.Range("AD27").Value = 1000
$ws.Range("L28").Value = 1.5
$ws.Range("M28").Value = 2.5
$ws.Range("G29").Value = 2.92
$ws.Range("H29").Value = 2.92
$ws.Range("I29").Value = 2.45
$ws.Range("L29").Value = 1.5
$ws.Range("M29").Value = 2.25
$ws.Range("N29").Value = 2.45
$ws.Range("O29").Value = 1.42
$ws.Range("P29").Value = 1.57
$ws.Range("Q29").Value = 2.12
$ws.Range("R29").Value = 2.07
$ws.Range("S29").Value = 1.6
$ws.Range("T29").Value = 6.8
$ws.Range("U29").Value = 13.5
$ws.Range("V29").Value = 11.5
$ws.Range("W29").Value = 37
$ws.Range("X29").Value = 32
$ws.Range("Y29").Value = 50
$ws.Range("Z29").Value = 6.3
$ws.Range("AA29").Value = 5.9
$ws.Range("AB29").Value = 19
$ws.Range("AC29").Value = 120
$ws.Range("AE29").Value = 5.9
$ws.Range("AF29").Value = 10.25
$ws.Range("AG29").Value = 10.25
$ws.Range("AH29").Value = 26
$ws.Range("AI29").Value = 26
$ws.Range("AJ29").Value = 50
$ws.Range("G30").Value = 2.37
$ws.Range("H30").Value = 3.25
$ws.Range("I30").Value = 2.77
$ws.Range("L30").Value = 1.31
$ws.Range("M30").Value = 2.87
$ws.Range("N30").Value = 1.91
$ws.Range("O30").Value = 1.7
$ws.Range("P30").Value = 1.42
$ws.Range("Q30").Value = 2.47
$ws.Range("R30").Value = 1.72
$ws.Range("S30").Value = 1.88
$ws.Range("T30").Value = 7.8
$ws.Range("U30").Value = 11.5
$ws.Range("V30").Value = 9.25
$ws.Range("W30").Value = 25
$ws.Range("X30").Value = 20
$ws.Range("Y30").Value = 30
$ws.Range("Z30").Value = 9.25
$ws.Range("AA30").Value = 6.3
$ws.Range("AB30").Value = 14
$ws.Range("AC30").Value = 70
$ws.Range("AD30").Value = 600
$ws.Range("AE30").Value = 8.5
$ws.Range("AF30").Value = 13.5
$ws.Range("AG30").Value = 10.25
$ws.Range("AH30").Value = 32
$ws.Range("AI30").Value = 24
$ws.Range("AJ30").Value = 35
